# Text editor Readme.docx edit:
# Insert a new paragraph "Default saving location är TextEditor/Bin/Debug."
# right after the "Grundläggande funktionaliteter..." paragraph and before
# the "Funktionaliteter:" paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that ends with "...och antas förstås som de är."
# (the "Grundläggande funktionaliteter..." paragraph) by searching for its
# distinctive trailing text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*och antas förstås som de är.*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Grundläggande funktionaliteter' paragraph"
}

# Insert a brand-new paragraph right after it.
$newPara = $target.Range.InsertParagraphAfter()

# Re-find the paragraph we just created: it's the one after $target now.
$insertedRange = $target.Next().Range

# Make sure the new paragraph carries the same language formatting as its
# neighbours (sv-SE), then set its text.
$insertedRange.Text = "Default saving location är TextEditor/Bin/Debug."
$insertedRange.LanguageID = 1053
